$wb = $excel.ActiveWorkbook
$dayWs = $wb.Worksheets.Item("day")
$weekWs = $wb.Worksheets.Item("week")

# --- Append new rows 249-259 to the "day" sheet ---
# Force column D to Text format for the new rows so the bsecode
# values are stored as strings (matching the scraped inlineStr data)
$dayWs.Range("D249:D259").NumberFormat = "@"

$dayWs.Cells.Item(249, 1).Value = 1
$dayWs.Cells.Item(249, 2).Value = "ABBOTINDIA"
$dayWs.Cells.Item(249, 3).Value = "Abbott India Limited"
$dayWs.Cells.Item(249, 4).Value = "500488"
$dayWs.Cells.Item(249, 5).Value = -0.9
$dayWs.Cells.Item(249, 6).Value = 27975
$dayWs.Cells.Item(249, 7).Value = 12019
$dayWs.Cells.Item(249, 8).Value = "day"
$dayWs.Cells.Item(249, 9).Value = "05/08/2024 11:34:38"

$dayWs.Cells.Item(250, 1).Value = 2
$dayWs.Cells.Item(250, 2).Value = "METROPOLIS"
$dayWs.Cells.Item(250, 3).Value = "Metropolis Healthcare Ltd"
$dayWs.Cells.Item(250, 4).Value = "542650"
$dayWs.Cells.Item(250, 5).Value = -1.44
$dayWs.Cells.Item(250, 6).Value = 2053.75
$dayWs.Cells.Item(250, 7).Value = 130042
$dayWs.Cells.Item(250, 8).Value = "day"
$dayWs.Cells.Item(250, 9).Value = "05/08/2024 11:34:38"

$dayWs.Cells.Item(251, 1).Value = 3
$dayWs.Cells.Item(251, 2).Value = "COROMANDEL"
$dayWs.Cells.Item(251, 3).Value = "Coromandel International Limited"
$dayWs.Cells.Item(251, 4).Value = "506395"
$dayWs.Cells.Item(251, 5).Value = -1.55
$dayWs.Cells.Item(251, 6).Value = 1613.6
$dayWs.Cells.Item(251, 7).Value = 417898
$dayWs.Cells.Item(251, 8).Value = "day"
$dayWs.Cells.Item(251, 9).Value = "05/08/2024 11:34:38"

$dayWs.Cells.Item(252, 1).Value = 4
$dayWs.Cells.Item(252, 2).Value = "CIPLA"
$dayWs.Cells.Item(252, 3).Value = "Cipla Limited"
$dayWs.Cells.Item(252, 4).Value = "500087"
$dayWs.Cells.Item(252, 5).Value = -0.9
$dayWs.Cells.Item(252, 6).Value = 1515.05
$dayWs.Cells.Item(252, 7).Value = 1719770
$dayWs.Cells.Item(252, 8).Value = "day"
$dayWs.Cells.Item(252, 9).Value = "05/08/2024 11:34:38"

$dayWs.Cells.Item(253, 1).Value = 5
$dayWs.Cells.Item(253, 2).Value = "VOLTAS"
$dayWs.Cells.Item(253, 3).Value = "Voltas Limited"
$dayWs.Cells.Item(253, 4).Value = "500575"
$dayWs.Cells.Item(253, 5).Value = -2.61
$dayWs.Cells.Item(253, 6).Value = 1451.85
$dayWs.Cells.Item(253, 7).Value = 1259451
$dayWs.Cells.Item(253, 8).Value = "day"
$dayWs.Cells.Item(253, 9).Value = "05/08/2024 11:34:38"

$dayWs.Cells.Item(254, 1).Value = 6
$dayWs.Cells.Item(254, 2).Value = "PVRINOX"
$dayWs.Cells.Item(254, 3).Value = "PVR Inox Ltd"
$dayWs.Cells.Item(254, 4).Value = "532689"
$dayWs.Cells.Item(254, 5).Value = -3.61
$dayWs.Cells.Item(254, 6).Value = 1417.8
$dayWs.Cells.Item(254, 7).Value = 416731
$dayWs.Cells.Item(254, 8).Value = "day"
$dayWs.Cells.Item(254, 9).Value = "05/08/2024 11:34:38"

$dayWs.Cells.Item(255, 1).Value = 7
$dayWs.Cells.Item(255, 2).Value = "ZYDUSLIFE"
$dayWs.Cells.Item(255, 3).Value = "Zydus Lifesciences Ltd"
$dayWs.Cells.Item(255, 4).Value = "532321"
$dayWs.Cells.Item(255, 5).Value = -1.79
$dayWs.Cells.Item(255, 6).Value = 1227.55
$dayWs.Cells.Item(255, 7).Value = 1115164
$dayWs.Cells.Item(255, 8).Value = "day"
$dayWs.Cells.Item(255, 9).Value = "05/08/2024 11:34:38"

$dayWs.Cells.Item(256, 1).Value = 8
$dayWs.Cells.Item(256, 2).Value = "RAMCOCEM"
$dayWs.Cells.Item(256, 3).Value = "The Ramco Cements Limited"
$dayWs.Cells.Item(256, 4).Value = "500260"
$dayWs.Cells.Item(256, 5).Value = -2.55
$dayWs.Cells.Item(256, 6).Value = 806.65
$dayWs.Cells.Item(256, 7).Value = 662526
$dayWs.Cells.Item(256, 8).Value = "day"
$dayWs.Cells.Item(256, 9).Value = "05/08/2024 11:34:38"

$dayWs.Cells.Item(257, 1).Value = 9
$dayWs.Cells.Item(257, 2).Value = "GUJGASLTD"
$dayWs.Cells.Item(257, 3).Value = "Gujarat Gas Limited"
$dayWs.Cells.Item(257, 4).Value = "539336"
$dayWs.Cells.Item(257, 5).Value = -2.27
$dayWs.Cells.Item(257, 6).Value = 642
$dayWs.Cells.Item(257, 7).Value = 861275
$dayWs.Cells.Item(257, 8).Value = "day"
$dayWs.Cells.Item(257, 9).Value = "05/08/2024 11:34:38"

$dayWs.Cells.Item(258, 1).Value = 10
$dayWs.Cells.Item(258, 2).Value = "GRANULES"
$dayWs.Cells.Item(258, 3).Value = "Granules India Limited"
$dayWs.Cells.Item(258, 4).Value = "532482"
$dayWs.Cells.Item(258, 5).Value = -0.68
$dayWs.Cells.Item(258, 6).Value = 641.7
$dayWs.Cells.Item(258, 7).Value = 3299558
$dayWs.Cells.Item(258, 8).Value = "day"
$dayWs.Cells.Item(258, 9).Value = "05/08/2024 11:34:38"

$dayWs.Cells.Item(259, 1).Value = 11
$dayWs.Cells.Item(259, 2).Value = "JUBLFOOD"
$dayWs.Cells.Item(259, 3).Value = "Jubilant Foodworks Limited"
$dayWs.Cells.Item(259, 4).Value = "533155"
$dayWs.Cells.Item(259, 5).Value = -3.51
$dayWs.Cells.Item(259, 6).Value = 585.4
$dayWs.Cells.Item(259, 7).Value = 2837100
$dayWs.Cells.Item(259, 8).Value = "day"
$dayWs.Cells.Item(259, 9).Value = "05/08/2024 11:34:38"

# --- Fix bsecode column type on "week" sheet rows 117-130: ---
# previously stored as text, should be numeric
$weekWs.Cells.Item(117, 4).Value = 500530
$weekWs.Cells.Item(118, 4).Value = 542652
$weekWs.Cells.Item(119, 4).Value = 532644
$weekWs.Cells.Item(120, 4).Value = 500480
$weekWs.Cells.Item(121, 4).Value = 511218
$weekWs.Cells.Item(122, 4).Value = 532454
$weekWs.Cells.Item(123, 4).Value = 524494
$weekWs.Cells.Item(124, 4).Value = 500260
$weekWs.Cells.Item(125, 4).Value = 540611
$weekWs.Cells.Item(126, 4).Value = 500850
$weekWs.Cells.Item(127, 4).Value = 500086
$weekWs.Cells.Item(128, 4).Value = 532720
$weekWs.Cells.Item(129, 4).Value = 540065
$weekWs.Cells.Item(130, 4).Value = 532822
